$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "61.020.50"
Set-TextValue "E2" "  +0.62%  "
Set-TextValue "D3" "2.664.32"
Set-TextValue "E3" "  +1.99%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "529.29"
Set-TextValue "E5" "  +3.07%  "
Set-TextValue "D6" "155.94"
Set-TextValue "E6" "  +0.93%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "D8" "0.581"
Set-TextValue "E8" "  -1.37%  "
Set-TextValue "E9" "  -2.95%  "
Set-TextValue "E10" "  +5.11%  "
Set-TextValue "D11" "0.352"
Set-TextValue "E11" "  +1.51%  "
Set-TextValue "E12" "  -0.43%  "
Set-TextValue "D13" "3.134.58"
Set-TextValue "E13" "  +2.02%  "
Set-TextValue "D14" "61.008.28"
Set-TextValue "E14" "  +0.65%  "
Set-TextValue "D15" "22.18"
Set-TextValue "E15" "  +2.34%  "
Set-TextValue "E16" "  +1.61%  "
Set-TextValue "D17" "2.666.36"
Set-TextValue "E17" "  +1.67%  "
Set-TextValue "E18" "  +0.99%  "
Set-TextValue "D19" "355.13"
Set-TextValue "E19" "  -0.71%  "
Set-TextValue "D20" "10.69"
Set-TextValue "E20" "  +0.64%  "
Set-TextValue "D21" "6.32"
Set-TextValue "E21" "  +2.26%  "
Set-TextValue "D22" "0.996"
Set-TextValue "E22" "  -0.24%  "
Set-TextValue "D23" "61.69"
Set-TextValue "E23" "  +1.53%  "
Set-TextValue "D24" "0.431"
Set-TextValue "E24" "  +1.80%  "
Set-TextValue "E25" "  +0.94%  "
Set-TextValue "D26" "0.998"
Set-TextValue "E26" "  +0.18%  "
Set-TextValue "D27" "0.0₃0857"
Set-TextValue "E27" "  +1.40%  "
Set-TextValue "E28" "  -0.28%  "
Set-TextValue "E29" "  -0.01%  "
Set-TextValue "D30" "6.18"
Set-TextValue "E30" "  +3.89%  "
Set-TextValue "D31" "19.54"
Set-TextValue "E31" "  +0.43%  "
Set-TextValue "D32" "1.62"
Set-TextValue "E32" "  +2.60%  "
Set-TextValue "D33" "150.08"
Set-TextValue "E33" "  -1.09%  "
Set-TextValue "E34" "  +3.16%  "
Set-TextValue "E35" "  +0.68%  "
Set-TextValue "D36" "0.923"
Set-TextValue "E36" "  +8.67%  "
Set-TextValue "D37" "0.896"
Set-TextValue "E37" "  +1.51%  "
Set-TextValue "D38" "36.91"
Set-TextValue "E38" "  +1.53%  "
Set-TextValue "D39" "306.00"
Set-TextValue "E39" "  +4.44%  "
Set-TextValue "E40" "  +0.24%  "
Set-TextValue "E41" "  +0.75%  "
Set-TextValue "D42" "0.645"
Set-TextValue "E42" "  +3.55%  "
Set-TextValue "E43" "  +0.54%  "
Set-TextValue "D44" "20.48"
Set-TextValue "E44" "  +3.21%  "
Set-TextValue "E45" "  +1.57%  "
Set-TextValue "E46" "  +0.05%  "
Set-TextValue "B47" "VeChain"
Set-TextValue "C47" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0243"
Set-TextValue "E47" "  +3.41%  "
Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "4.98"
Set-TextValue "E48" "  +0.95%  "
Set-TextValue "E49" "  +6.95%  "
Set-TextValue "D50" "10.37"
Set-TextValue "E50" "  +0.61%  "
Set-TextValue "D51" "2.005.77"
Set-TextValue "E51" "  +0.55%  "
